# Apply the "Ilkal removed from SuppTable1" edit:
#  - delete the row holding the "SUSCEPTIBILITY OF 4 SPECIES OF MOSQUITOS TO
#    CHANDIPURA VIRUS AND ITS DETECTION BY IMMUNOFLUORESCENCE" title (row 173,
#    which has no DOI in column B)
#  - rename the worksheet from "Emily2" to "SuppTable1"
#  - move the active selection from E7 to A7
#  - re-apply the existing sort over A2:B(last row) so the sheet's stored
#    sortState/sortCondition range shrinks along with the data

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Locate the row to remove by matching its title text, rather than hard-coding
# the row index, so the script is robust to any other differences.
$targetText = "SUSCEPTIBILITY OF 4 SPECIES OF MOSQUITOS TO CHANDIPURA VIRUS AND ITS DETECTION BY IMMUNOFLUORESCENCE"
$found = $ws.Columns.Item(1).Find($targetText)
if ($found -ne $null) {
    $found.EntireRow.Delete()
}

# Rename the sheet.
$ws.Name = "SuppTable1"

# Update the visible selection.
$ws.Range("A7").Select()

# Recompute the used range on column A/B and re-apply the sort so the stored
# sortState reference shrinks from B267 to B266.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$sortRange = $ws.Range("A2:B" + $lastRow)
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("A2:A" + $lastRow))
$ws.Sort.SetRange($sortRange)
$ws.Sort.Header = 2
$ws.Sort.Apply()
